# Generate Report for Handoff
# Update the localization status report:
#   - Status changes from "In Translation" to "Ready for handoff"
#   - Handoff timestamps are refreshed to reflect the new report run

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-21 14:36:10"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-21 14:36:07"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-21 14:36:10"
